$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (de)
$ws.Range("B2").Value = 17
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "xl_ent"

# Update existing row 3 (en)
$ws.Range("B3").Value = 307
$ws.Range("C3").Value = 25
$ws.Range("D3").Value = "abitofzen"

# Update existing row 4 (es)
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 0

# Update existing row 5 (fr)
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "scoub"

# Insert new row for "it" before current row 6 ("ja"), shifting rows down
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "it"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = ""

# Current row 7 is "ja" (originally row 6) - update values
$ws.Range("A7").Value = "ja"
$ws.Range("B7").Value = 31
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "ichigaya2016"

# Insert new row for "ko" before current row 8 ("nl"), shifting rows down
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "ko"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = ""

# Current row 9 is "nl" (originally row 7)
$ws.Range("A9").Value = "nl"
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "NuOpRadio2"

# Current row 10 is "pt" (originally row 8)
$ws.Range("A10").Value = "pt"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = ""

# Insert new row for "qme" before current row 11 ("ru"), shifting rows down
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "qme"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = ""

# Current row 12 is "ru" (originally row 9)
$ws.Range("A12").Value = "ru"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = ""

# Current row 13 is "tr" (originally row 10)
$ws.Range("A13").Value = "tr"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = ""

# Insert new row for "uk" before current row 14 ("zh"), shifting rows down
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "uk"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = ""

# Current row 15 is "zh" (originally row 11)
$ws.Range("A15").Value = "zh"
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = ""

# Current row 16 is "zh-CN" (originally row 12)
$ws.Range("A16").Value = "zh-CN"
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = ""

# Current row 17 is "Total" (originally row 13)
$ws.Range("A17").Value = "Total"
$ws.Range("B17").Value = 413
$ws.Range("C17").Value = 33
$ws.Range("D17").Value = ""
